# Generate Report for Handoff
#
# The localization file "0e314636-8642-429d-95e2-56fccc4a9f14.md" was
# handed off again (status moved from "Handed back: in sync with en-US"
# to "Ready for handoff" with new timestamps), while
# "5f2c9c7b-8ee7-48fd-af0c-5866b81aa82d.md" stays "Handed back" as-is.
# The status report table is regenerated, which re-orders the two data
# rows on every sheet: row 2 now carries the 5f2c9c7b file, row 3 now
# carries the 0e314636 file (with its refreshed status/timestamps).

$wb = $excel.ActiveWorkbook

$ZH_MD_0E   = "0e314636-8642-429d-95e2-56fccc4a9f14.md"
$ZH_MD_5F   = "5f2c9c7b-8ee7-48fd-af0c-5866b81aa82d.md"
$STATUS_BACK = "Handed back: in sync with en-US"
$STATUS_READY = "Ready for handoff"

# ---------------------------------------------------------------
# Sheet: Overview
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A2").Value = $ZH_MD_5F
$ws.Range("B2").Value = $STATUS_BACK
$ws.Range("C2").Value = $STATUS_BACK
$ws.Range("D2").Value = "2016-03-22 08:45:14"

$ws.Range("A3").Value = $ZH_MD_0E
$ws.Range("B3").Value = $STATUS_READY
$ws.Range("C3").Value = $STATUS_READY
$ws.Range("D3").Value = "2016-03-22 08:46:41"

foreach ($hl in $ws.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') {
        $hl.TextToDisplay = $ZH_MD_5F
    } elseif ($addr -eq '$A$3') {
        $hl.TextToDisplay = $ZH_MD_0E
    }
}

# ---------------------------------------------------------------
# Sheet: zh-cn
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ZH_XLF_0E = "0e314636-8642-429d-95e2-56fccc4a9f14.2981ce20929d003ce22b02035c8278eea0ddbf86.zh-cn.xlf"
$ZH_XLF_5F = "5f2c9c7b-8ee7-48fd-af0c-5866b81aa82d.09e02514cb63c024c396759722709106f98eec03.zh-cn.xlf"

$ws.Range("A2").Value = $ZH_MD_5F
$ws.Range("B2").Value = ".md"
$ws.Range("C2").Value = $STATUS_BACK
$ws.Range("D2").Value = $ZH_XLF_5F
$ws.Range("E2").Value = "2016-03-22 08:45:10"
$ws.Range("F2").Value = $ZH_MD_5F
$ws.Range("G2").Value = $ZH_XLF_5F
$ws.Range("H2").Value = "2016-03-22 08:45:49"
$ws.Range("J2").Value = "Include"

$ws.Range("A3").Value = $ZH_MD_0E
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = $STATUS_READY
$ws.Range("D3").Value = $ZH_XLF_0E
$ws.Range("E3").Value = "2016-03-22 08:46:37"
$ws.Range("F3").Value = $ZH_MD_0E
$ws.Range("G3").Value = $ZH_XLF_0E
$ws.Range("H3").Value = "2016-03-22 08:45:49"
$ws.Range("J3").Value = "Include"

foreach ($hl in $ws.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') {
        $hl.TextToDisplay = $ZH_MD_5F
    } elseif ($addr -eq '$D$2') {
        $hl.TextToDisplay = $ZH_XLF_5F
    } elseif ($addr -eq '$F$2') {
        $hl.TextToDisplay = $ZH_MD_5F
    } elseif ($addr -eq '$G$2') {
        $hl.TextToDisplay = $ZH_XLF_5F
    } elseif ($addr -eq '$A$3') {
        $hl.TextToDisplay = $ZH_MD_0E
    } elseif ($addr -eq '$D$3') {
        $hl.TextToDisplay = $ZH_XLF_0E
    } elseif ($addr -eq '$F$3') {
        $hl.TextToDisplay = $ZH_MD_0E
    } elseif ($addr -eq '$G$3') {
        $hl.TextToDisplay = $ZH_XLF_0E
    }
}

# ---------------------------------------------------------------
# Sheet: de-de
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$DE_XLF_0E = "0e314636-8642-429d-95e2-56fccc4a9f14.2981ce20929d003ce22b02035c8278eea0ddbf86.de-de.xlf"
$DE_XLF_5F = "5f2c9c7b-8ee7-48fd-af0c-5866b81aa82d.09e02514cb63c024c396759722709106f98eec03.de-de.xlf"

$ws.Range("A2").Value = $ZH_MD_5F
$ws.Range("B2").Value = ".md"
$ws.Range("C2").Value = $STATUS_BACK
$ws.Range("D2").Value = $DE_XLF_5F
$ws.Range("E2").Value = "2016-03-22 08:45:14"
$ws.Range("F2").Value = $ZH_MD_5F
$ws.Range("G2").Value = $DE_XLF_5F
$ws.Range("H2").Value = "2016-03-22 08:45:57"
$ws.Range("J2").Value = "Include"

$ws.Range("A3").Value = $ZH_MD_0E
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = $STATUS_READY
$ws.Range("D3").Value = $DE_XLF_0E
$ws.Range("E3").Value = "2016-03-22 08:46:41"
$ws.Range("F3").Value = $ZH_MD_0E
$ws.Range("G3").Value = $DE_XLF_0E
$ws.Range("H3").Value = "2016-03-22 08:45:57"
$ws.Range("J3").Value = "Include"

foreach ($hl in $ws.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') {
        $hl.TextToDisplay = $ZH_MD_5F
    } elseif ($addr -eq '$D$2') {
        $hl.TextToDisplay = $DE_XLF_5F
    } elseif ($addr -eq '$F$2') {
        $hl.TextToDisplay = $ZH_MD_5F
    } elseif ($addr -eq '$G$2') {
        $hl.TextToDisplay = $DE_XLF_5F
    } elseif ($addr -eq '$A$3') {
        $hl.TextToDisplay = $ZH_MD_0E
    } elseif ($addr -eq '$D$3') {
        $hl.TextToDisplay = $DE_XLF_0E
    } elseif ($addr -eq '$F$3') {
        $hl.TextToDisplay = $ZH_MD_0E
    } elseif ($addr -eq '$G$3') {
        $hl.TextToDisplay = $DE_XLF_0E
    }
}
